$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MSRP (column D) and DPHF (column E) values for LX / GX rows
# Row 29 - GX 460 (9700)
$ws.Range("D29").Value = 53100

# Row 30 - GX 460 Premium (9700PM)
$ws.Range("D30").Value = 55890

# Row 31 - LX 570 Two-Row (9710)
$ws.Range("D31").Value = 64365

# Row 32 - LX 570 Three-Row (9625)
$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025

# Row 33 - LX 570 Inspiration Series SE base (9620)
$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

# Row 34 - LX 570 Inspiration Series SE (9620 (SE)) - activate MSRP value
$ws.Range("D34").Value = 99310
$ws.Range("D34").NumberFormat = $ws.Range("D33").NumberFormat
$ws.Range("E34").Value = 1025

# Row 43 shrink from 48pt to 45pt
$ws.Rows.Item(43).RowHeight = 45

# Update the active view state to match the latest edits
$ws.Range("D29").Select()
